# Testcases.xlsx edit script
# Implements: status column (Enabled/Disabled), enabling all 3 testcases on
# Sheet1, and switching the hard-coded WordPress URL for a global "testUrl"
# variable; also refreshes / extends the underlying test-step text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1: collapse the old 6-row / 2-col layout (TC_001 plus five loose
# scratch rows) into a clean 3-row / 3-col table: Status | TestCase | Steps
# ---------------------------------------------------------------------

# Drop the three now-unused scratch rows (old rows 4-6) first, from the
# bottom up so row indices of earlier rows stay stable.
$ws1.Rows.Item(6).Delete()
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

# Insert a new first column to hold the Enabled/Disabled status.
$ws1.Columns.Item(1).Insert()

# Status values - set in this specific order (TC_002's "Enabled" first) so
# the workbook's shared-string table is rebuilt in the same order the
# original author produced it in.
$ws1.Range("A2").Value = "Enabled"
$ws1.Range("A1").Value = "Disabled"
$ws1.Range("A3").Value = "Disabled"

# TestCase id column
$ws1.Range("B1").Value = "TC_001"
$ws1.Range("B2").Value = "TC_002"
$ws1.Range("B3").Value = "TC_003"

# Row 1: TC_001 steps - switched from a hard-coded URL to the global
# "testUrl" config value.
$ws1.Range("C1").Value = "Invoke browser`nLoad testUrl`nType admin in UserNameField with id = user_login`nClear UserNameField with id = user_login`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nGet tagname of LogInBtn with id = wp-submit Assign {LogInBtnTagName}`n{LogInBtnTagName} VerifyEqual input`nGet value of attribute: value of LogInBtn with id = wp-submit  Assign {LogInBtnValue}`n{LogInBtnValue} VerifyEqual Log In`nIs LogInBtn with id = wp-submit displayed Assign {LogInBtnDisplayed}`n{LogInBtnDisplayed} VerifyEqual true`nIs LogInBtn with id = wp-submit enabled Assign {LogInBtnEnabled}`n{LogInBtnEnabled} VerifyEqual true`nGet text of UserNameLabel with xpath = //form[@id='loginform']/p[1]/label Assign {UserNameLabelText}`n{UserNameLabelText} VerifyEqual  Username`nIs RememberMeChkBox with id = rememberme selected Assign {RememberMeSelected}`n{RememberMeSelected} VerifyEqual false`nClick RememberMeChkBox with id = rememberme`nIs RememberMeChkBox with id = rememberme selected Assign {RememberMeSelected}`n{RememberMeSelected} VerifyEqual true`nClick LogInBtn with id = wp-submit`nGet page title Assign {title}`n{title} VerifyEqual Dashboard ‹ test — WordPress`nGet current page url Assign {pageUrl}`n{pageUrl} VerifyEqual http://127.0.0.1/wordpress/wp-admin/`nQuit browser"

# Row 2: TC_002 steps - new publish-a-post scenario
$ws1.Range("C2").Value = "Invoke browser`nLoad testUrl`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nClick LogInBtn with id = wp-submit`nType test post title in PostTitleField with id = title`nType test post content in PostContentField with id = content`nClick PublishBtn with id = publish`nSleep 2`nIs ViewPostLink with linkText = View post displayed Assign {ViewPostLinkDisplayed}`n{ViewPostLinkDisplayed} VerifyEqual true`nIs EditPostLink with linkText = Edit post displayed Assign {EditPostLinkDisplayed}`n{EditPostLinkDisplayed} VerifyEqual true`nQuit browser"

# Row 3: TC_003 steps - new edit/trash-a-post scenario (trailing blank
# line preserved exactly as authored)
$ws1.Range("C3").Value = "Invoke browser`nLoad testUrl`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nClick LogInBtn with id = wp-submit`nClick PostsLink with xpath = //li[@id='menu-posts']/a/div[3]`nGet page title Assign {title}`n{title} VerifyEqual Posts ‹ test — WordPress`nSelect Edit from ActionsListBox with xpath = //form[@id='posts-filter']/div[1]/div[1]/select`nClick AddedPostLink with xpath = //tr[contains(@id,'post')]/td[1]/strong/a`nGet page title Assign {title}`n{title} VerifyEqual Edit Post ‹ test — WordPress`nClick MoveToTrashLink with xpath = //div[@id='delete-action']/a`nQuit browser`n"

# Formatting for the new status column: centered text, no wrap. Set the
# alignment on A1 only, then fan it out to A2:A3 via a formats-only paste so
# the style table doesn't pick up intermediate (H-only / V-only) entries.
$ws1.Range("A1").HorizontalAlignment = -4108
$ws1.Range("A1").VerticalAlignment = -4108
$ws1.Range("A1").Copy()
$ws1.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Columns.Item(1).ColumnWidth = 12

# Row heights
$ws1.Rows.Item(1).RowHeight = 409.5
$ws1.Rows.Item(2).RowHeight = 210
$ws1.Rows.Item(3).RowHeight = 225

# Selection / view state
$ws1.Range("C2").Select()

# ---------------------------------------------------------------------
# Sheet2: no content change (shared-string renumbering only, handled
# automatically by the engine's own string table on save).
# ---------------------------------------------------------------------
